# Edit script: insert 3 new data rows (861-863) into the weekly price table,
# pushing the previously-existing rows 861-927 down to 864-930.
#
# The workbook has a single worksheet with a data table starting at row 2
# (row 1 looks like a header/first-record row depending on export) and the
# rows we care about are 861..927 before the edit (861..930 after).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 861. This shifts the former rows
# 861-927 down to 864-930, preserving all of their existing values/styles.
$ws.Rows("861:863").Insert()

# Now populate the 3 newly inserted rows with their final values.

# Row 861
$ws.Cells.Item(861, 1).Value = 10
$ws.Cells.Item(861, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(861, 3).Value = "La Araucanía"
$ws.Cells.Item(861, 4).Value = 45166
$ws.Cells.Item(861, 5).Value = 9
$ws.Cells.Item(861, 6).Value = 100112032
$ws.Cells.Item(861, 7).Value = "Zapallo italiano"
$ws.Cells.Item(861, 8).Value = "Bola 8"
$ws.Cells.Item(861, 9).Value = "Primera"
$ws.Cells.Item(861, 10).Value = 50
$ws.Cells.Item(861, 11).Value = 18000
$ws.Cells.Item(861, 12).Value = 18000
$ws.Cells.Item(861, 13).Value = 18000
$ws.Cells.Item(861, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(861, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(861, 16).Value = 360
$ws.Cells.Item(861, 17).Value = 50
$ws.Cells.Item(861, 18).Value = "Hortaliza"

# Row 862
$ws.Cells.Item(862, 1).Value = 10
$ws.Cells.Item(862, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(862, 3).Value = "La Araucanía"
$ws.Cells.Item(862, 4).Value = 45166
$ws.Cells.Item(862, 5).Value = 9
$ws.Cells.Item(862, 6).Value = 100112032
$ws.Cells.Item(862, 7).Value = "Zapallo italiano"
$ws.Cells.Item(862, 8).Value = "Huracán"
$ws.Cells.Item(862, 9).Value = "Primera"
$ws.Cells.Item(862, 10).Value = 80
$ws.Cells.Item(862, 11).Value = 18000
$ws.Cells.Item(862, 12).Value = 18000
$ws.Cells.Item(862, 13).Value = 18000
$ws.Cells.Item(862, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(862, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(862, 16).Value = 360
$ws.Cells.Item(862, 17).Value = 50
$ws.Cells.Item(862, 18).Value = "Hortaliza"

# Row 863
$ws.Cells.Item(863, 1).Value = 10
$ws.Cells.Item(863, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(863, 3).Value = "La Araucanía"
$ws.Cells.Item(863, 4).Value = 45166
$ws.Cells.Item(863, 5).Value = 9
$ws.Cells.Item(863, 6).Value = 100112032
$ws.Cells.Item(863, 7).Value = "Zapallo italiano"
$ws.Cells.Item(863, 8).Value = "Sin especificar"
$ws.Cells.Item(863, 9).Value = "Primera"
$ws.Cells.Item(863, 10).Value = 550
$ws.Cells.Item(863, 11).Value = 15000
$ws.Cells.Item(863, 12).Value = 17000
$ws.Cells.Item(863, 13).Value = 15545
$ws.Cells.Item(863, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(863, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(863, 16).Value = 311
$ws.Cells.Item(863, 17).Value = 50
$ws.Cells.Item(863, 18).Value = "Hortaliza"
